# data provider added in test
# Update the "Tweets" worksheet (Sheet2) text values, widen column A,
# and move the active selection from A5 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Make sure Sheet2 is the active/selected sheet (it was already the
# active tab in the source workbook).
[void]$ws.Activate()

# Update the shared-string text used by the data-provider rows.
$ws.Range("A2").Value = "Hello All Indians"
$ws.Range("A3").Value = "Good Evening India and Noida"
$ws.Range("A4").Value = "Automation Framework TDD"

# Widen column A to fit the new, longer text.
$ws.Columns.Item(1).ColumnWidth = 25.17

# Move the stored selection / active cell from A5 to A4.
[void]$ws.Range("A4").Select()
